# Generate Report for handoff
# - Update status of 2ac41cb1-... file from "Handed back: in sync with en-US" to "Ready for handoff"
#   and refresh its "Latest Handoff Datetime" on the zh-cn / de-de sheets.
# - Remove the row belonging to the 340eefdd-... file (it is no longer tracked), on all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# Row 3 belongs to 340eefdd-...; delete it so the .localization-config row shifts up.
$ws1.Rows.Item(3).Delete()

$ws1.Hyperlinks.Delete()
$h1 = $ws1.Hyperlinks
$l = $h1.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
$l = $h1.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config")
$l.TextToDisplay = ".localization-config"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-01-25 03:40:15"

# Row 3 belongs to 340eefdd-...; delete it so the .localization-config row shifts up.
$ws2.Rows.Item(3).Delete()

$ws2.Hyperlinks.Delete()
$h2 = $ws2.Hyperlinks
$l = $h2.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
$l = $h2.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/547e5119a31a5aa0c372863ce57cca36cea6165a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf"
$l = $h2.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d8bb2cf177853f6f770e7db892574102310e017b/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
$l = $h2.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/721f89fb0de81fd69cbc002c3d9c4657112c6a87/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.zh-cn.xlf"
$l = $h2.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config")
$l.TextToDisplay = ".localization-config"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-01-25 03:40:25"

# Row 3 belongs to 340eefdd-...; delete it so the .localization-config row shifts up.
$ws3.Rows.Item(3).Delete()

$ws3.Hyperlinks.Delete()
$h3 = $ws3.Hyperlinks
$l = $h3.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
$l = $h3.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/414197c5ac1b0ab62b3841ee1b86f4ee1b7e736a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf"
$l = $h3.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/cbed1a67d5766a9360015ca5a30faa1ab11215cf/e2e/2ac41cb1-a240-442c-a524-4efa10ef2ca6.md")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.md"
$l = $h3.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ed00aa0e9ea8db3a5b9f4ddfe2f53d419935d640/ol-handback/OpenLocalizationTestOrg/oltest.de-de/qimu/2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf")
$l.TextToDisplay = "2ac41cb1-a240-442c-a524-4efa10ef2ca6.df4bfcc7f35b9f2688b33325409ea6aadb12d18f.de-de.xlf"
$l = $h3.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b79fd7963f086c2aabea7f6f405ee180f4d2f276/.localization-config")
$l.TextToDisplay = ".localization-config"

Write-Output "Report regenerated for handoff"
